{"js": "// Apply the \"Something in the Closet\" concept-doc rework to the\n// \"Resources\" (to-do) bullet list at the end of the document.\n//\n// Before:\n//   - Implement player HUD.\n//     - Pause menu                                  [bookmark _GoBack]\n//   - User being able to move under the covers and peak.\n//   - Cinematic implementation.\n//   - Game over state implementation.\n//   - Update models and art.\n//\n// After:\n//   - User being able to move under the covers and peak.\n//   - Make entering the bed a smoother transition\n//     - Done the sequencer now just needs adjustment. [bookmark _GoBack]\n//   - Add under cover visual effect.\n//   - Cinematic implementation.\n//   - Game over state implementation.//Rethink this\n//     - No game over, think more PT\n//     - Negative effects to being seen by the monster\n//     - Negative effects to being caught by the monster\n//     - Maybe special ending for looking at it X amount of time without being seen.\n//   - Update models and art.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the six bullets that anchor this edit by their (unique) text,\n// rather than hard-coded indices, so the script is resilient to minor\n// paragraph-count differences elsewhere in the document.\nconst texts = paragraphs.items.map((p) => p.text);\n\nfunction findIndex(expected) {\n  const idx = texts.indexOf(expected);\n  if (idx === -1) {\n    throw new Error(\"Could not find paragraph with text: \" + expected);\n  }\n  return idx;\n}\n\nconst idxHud = findIndex(\"Implement player HUD. \");\nconst idxPause = findIndex(\"Pause menu\");\nconst idxCovers = findIndex(\"User being able to move under the covers and peak. \");\nconst idxGameOver = findIndex(\"Game over state implementation. \");\n\n// 1) \"Implement player HUD.\" -> \"User being able to move under the covers and peak.\"\nparagraphs.items[idxHud]\n  .getRange(Word.RangeLocation.whole)\n  .insertText(\"User being able to move under the covers and peak. \", Word.InsertLocation.replace);\nawait context.sync();\n\n// 2) New top-level bullet right after it: \"Make entering the bed a smoother transition\"\nparagraphs.items[idxHud].insertParagraph(\n  \"Make entering the bed a smoother transition\",\n  Word.InsertLocation.after\n);\nawait context.sync();\n\n// 3) \"Pause menu\" -> \"Done the sequencer now just needs adjustment.\" (keeps its bookmark)\nparagraphs.items[idxPause]\n  .getRange(Word.RangeLocation.whole)\n  .insertText(\"Done the sequencer now just needs adjustment.\", Word.InsertLocation.replace);\nawait context.sync();\n\n// 4) The old \"User being able to move under the covers and peak.\" bullet becomes\n//    \"Add under cover visual effect.\"\nparagraphs.items[idxCovers]\n  .getRange(Word.RangeLocation.whole)\n  .insertText(\"Add under cover visual effect.\", Word.InsertLocation.replace);\nawait context.sync();\n\n// 5) \"Game over state implementation. \" gains a second run \"//Rethink this\"\n//    (two separate <w:r> elements within the same paragraph).\n{\n  const p = paragraphs.items[idxGameOver];\n  const whole = p.getRange(Word.RangeLocation.whole);\n  const ooxml =\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData>' +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    '<w:body>' +\n    '<w:p>' +\n    '<w:pPr><w:pStyle w:val=\"ListParagraph\"/><w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"1\"/></w:numPr></w:pPr>' +\n    '<w:r><w:t xml:space=\"preserve\">Game over state implementation. </w:t></w:r>' +\n    '<w:r><w:t>//Rethink this</w:t></w:r>' +\n    '</w:p>' +\n    '</w:body></w:document>' +\n    '</pkg:xmlData></pkg:part></pkg:package>';\n  whole.insertOoxml(ooxml, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 6) Four new second-level bullets after the \"Game over state implementation.\" item.\nparagraphs.load(\"items/text\");\nawait context.sync();\nconst idxGameOver2 = paragraphs.items\n  .map((p) => p.text)\n  .indexOf(\"Game over state implementation. //Rethink this\");\n\nconst newSubBullets = [\n  \"No game over, think more PT\",\n  \"Negative effects to being seen by the monster\",\n  \"Negative effects to being caught by the monster\",\n  \"Maybe special ending for looking at it X amount of time without being seen.\",\n];\n\nlet anchor = paragraphs.items[idxGameOver2];\nfor (const bulletText of newSubBullets) {\n  const newPara = anchor.insertParagraph(bulletText, Word.InsertLocation.after);\n  await context.sync();\n\n  // Re-fetch the freshly-inserted paragraph (stale references don't report\n  // ranges/list info reliably) and bump it to the second outline level.\n  paragraphs.load(\"items/text\");\n  await context.sync();\n  const freshIdx = paragraphs.items.map((p) => p.text).indexOf(bulletText);\n  const freshPara = paragraphs.items[freshIdx];\n  freshPara.listItemOrNullObject.level = 1;\n  await context.sync();\n\n  anchor = freshPara;\n}\n", "ps1": "# Apply the \"Something in the Closet\" concept-doc rework to the\n# \"Resources\" (to-do) bullet list at the end of the document.\n#\n# Before:\n#   - Implement player HUD.\n#     - Pause menu                                  [bookmark _GoBack]\n#   - User being able to move under the covers and peak.\n#   - Cinematic implementation.\n#   - Game over state implementation.\n#   - Update models and art.\n#\n# After:\n#   - User being able to move under the covers and peak.\n#   - Make entering the bed a smoother transition\n#     - Done the sequencer now just needs adjustment. [bookmark _GoBack]\n#   - Add under cover visual effect.\n#   - Cinematic implementation.\n#   - Game over state implementation.//Rethink this\n#     - No game over, think more PT\n#     - Negative effects to being seen by the monster\n#     - Negative effects to being caught by the monster\n#     - Maybe special ending for looking at it X amount of time without being seen.\n#   - Update models and art.\n\n$d = $word.ActiveDocument\n\nfunction Find-ParagraphIndex($doc, $text) {\n    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {\n        $paraText = $doc.Paragraphs.Item($i).Range.Text.TrimEnd([char]13, [char]7)\n        if ($paraText -eq $text) {\n            return $i\n        }\n    }\n    throw \"Could not find paragraph with text: $text\"\n}\n\n# Resolve every index we need BEFORE mutating anything: step 1 renames\n# \"Implement player HUD.\" to text that is identical to the (still separate)\n# \"User being able...\" bullet below it, so searching by text after that\n# rename would be ambiguous. Grabbing all anchors up front avoids that.\n$idxHud = Find-ParagraphIndex $d \"Implement player HUD. \"\n$idxPause = Find-ParagraphIndex $d \"Pause menu\"\n$idxCovers = Find-ParagraphIndex $d \"User being able to move under the covers and peak. \"\n$idxGameOver = Find-ParagraphIndex $d \"Game over state implementation. \"\n\n# 1) \"Implement player HUD.\" -> \"User being able to move under the covers and peak.\"\n$d.Paragraphs.Item($idxHud).Range.Text = \"User being able to move under the covers and peak. \"\n\n# 2) New top-level bullet right after it: \"Make entering the bed a smoother transition\"\n$d.Paragraphs.Item($idxHud).Range.InsertParagraphAfter()\n$d.Paragraphs.Item($idxHud + 1).Range.Text = \"Make entering the bed a smoother transition\"\n\n# Inserting that extra paragraph shifted every subsequent paragraph index down by one.\n$idxPause = $idxPause + 1\n$idxCovers = $idxCovers + 1\n$idxGameOver = $idxGameOver + 1\n\n# 3) \"Pause menu\" -> \"Done the sequencer now just needs adjustment.\" (keeps its bookmark)\n$d.Paragraphs.Item($idxPause).Range.Text = \"Done the sequencer now just needs adjustment.\"\n\n# 4) The old \"User being able to move under the covers and peak.\" bullet becomes\n#    \"Add under cover visual effect.\"\n$d.Paragraphs.Item($idxCovers).Range.Text = \"Add under cover visual effect.\"\n\n# 5) \"Game over state implementation. \" gains a second run \"//Rethink this\"\n#    (two separate runs within the same paragraph) via a raw-OOXML replace of\n#    the paragraph's whole range (this keeps pPr/list formatting intact).\n$gameOverRange = $d.Paragraphs.Item($idxGameOver).Range\n$gameOverPkg = '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData>' +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    '<w:body>' +\n    '<w:p>' +\n    '<w:pPr><w:pStyle w:val=\"ListParagraph\"/><w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"1\"/></w:numPr></w:pPr>' +\n    '<w:r><w:t xml:space=\"preserve\">Game over state implementation. </w:t></w:r>' +\n    '<w:r><w:t>//Rethink this</w:t></w:r>' +\n    '</w:p>' +\n    '</w:body></w:document>' +\n    '</pkg:xmlData></pkg:part></pkg:package>'\n$gameOverRange.InsertXML($gameOverPkg)\n\n# 6) Four new second-level bullets after the \"Game over state implementation.\" item.\n#    InsertXML replaced that paragraph's content in place, so its paragraph\n#    index ($idxGameOver) is still valid.\n$newSubBullets = @(\n    \"No game over, think more PT\",\n    \"Negative effects to being seen by the monster\",\n    \"Negative effects to being caught by the monster\",\n    \"Maybe special ending for looking at it X amount of time without being seen.\"\n)\n\n$anchorIndex = $idxGameOver\nforeach ($bulletText in $newSubBullets) {\n    $d.Paragraphs.Item($anchorIndex).Range.InsertParagraphAfter()\n    $newIndex = $anchorIndex + 1\n    $newParagraph = $d.Paragraphs.Item($newIndex)\n    $newParagraph.Range.Text = $bulletText\n    $newParagraph.Range.ListFormat.ListLevelNumber = 2\n    $anchorIndex = $newIndex\n}\n"}
